$wb = $excel.ActiveWorkbook

# --- Text change: "Ready for handoff" -> "In Translation" ---
# Overview sheet: per-language status columns (E = zh-cn, F = de-de)
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2:F4").Value = "In Translation"

# zh-cn sheet: Status column (C)
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2:C4").Value = "In Translation"

# de-de sheet: Status column (C)
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2:C4").Value = "In Translation"

# --- Column width change following the shorter status text ---
# The Status column(s) narrow from ~17.22 to ~13.41 "characters" stored
# width. This engine quantizes Range.ColumnWidth to the nearest 1/6th of a
# character before persisting it (stored = round_1_6(ColumnWidth) + 5/6),
# so 12.5 is the input that lands closest on the saved grid to 13.4101845877511.
$newStatusWidth = 12.5

$wsOverview.Columns.Item(5).ColumnWidth = $newStatusWidth
$wsOverview.Columns.Item(6).ColumnWidth = $newStatusWidth

$wsZhCn.Columns.Item(3).ColumnWidth = $newStatusWidth

$wsDeDe.Columns.Item(3).ColumnWidth = $newStatusWidth
